$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 37

# Headline text (identical content for daily_headlines and daily_headlines_clean)
$headline = @'
invesco ltd acquires share of archrock inc nysearoc invesco ltd lower position in datadog inc nasdaqddog invesco ltd sell share of enbridge inc nyseenb invesco ltd lower stake in datadog inc nasdaqddog palo alto network inc nasdaqpanw share sold by entropy technology lp entropy technology lp sell share of palo alto network inc nasdaqpanw entropy technology lp raise stock holding in plain gp holding lp nysepagp entropy technology lp purchase new share in snap inc nysesnap entropy technology lp purchase share of makemytrip limited nasdaqmmyt entropy technology lp make new investment in coinbase global inc nasdaqcoin entropy technology lp make new investment in spx technology inc nysespxc entropy technology lp invests in makemytrip limited nasdaqmmyt entropy technology lp take position in innodata inc nasdaqinod entropy technology lp sell share of mohawk industry inc nysemhk share in coinbase global inc nasdaqcoin purchased by entropy technology lp entropy technology lp buy new share in innodata inc nasdaqinod share in glacier bancorp inc nasdaqgbci purchased by entropy technology lp entropy technology lp acquires share of biomarin pharmaceutical inc nasdaqbmrn entropy technology lp purchase new stake in integer holding corporation nyseitgr share in glacier bancorp inc nasdaqgbci bought by entropy technology lp share in cellebrite di ltd nasdaqclbt bought by entropy technology lp reason pfizers yielding dividend is getting safer possible reason warren buffett shunned his favorite stock for the fourth straight quarter despite sitting on billion in cash mitsubishi ufj asset management co ltd boost stock position in urban edge property nyseue li auto inc sponsored adr nasdaqli share sold by mitsubishi ufj asset management co ltd swiss national bank ha million stock holding in interpublic group of company inc the nyseipg mitsubishi ufj asset management co ltd purchase share of bwx technology inc nysebwxt medical property trust inc nysempw share purchased by mitsubishi ufj asset management co ltd mitsubishi ufj asset management co ltd raise holding in urban edge property nyseue news corporation nasdaqnwsa share sold by swiss national bank outfront medium inc nyseout share acquired by mitsubishi ufj asset management co ltd mitsubishi ufj asset management co ltd ha million position in sl green realty corporation nyseslg mitsubishi ufj asset management co ltd decrease stake in sprout farmer market inc nasdaqsfm mitsubishi ufj asset management co ltd buy share of millrose property inc nysemrp mitsubishi ufj asset management co ltd cut holding in lucid group inc nasdaqlcid mitsubishi ufj asset management co ltd reduces stock position in liquidia technology inc nasdaqlqda mitsubishi ufj asset management co ltd ha million stock position in sprout farmer market inc nasdaqsfm mitsubishi ufj asset management co ltd increase stock holding in ducommun incorporated nysedco mitsubishi ufj asset management co ltd reduces stake in lucid group inc nasdaqlcid li auto inc sponsored adr nasdaqli share sold by mitsubishi ufj asset management co ltd mitsubishi ufj asset management co ltd ha million stake in sabra healthcare reit inc nasdaqsbra china target lithuanian bank in response to new eu sanction procept biorobotics corporation nasdaqprct share bought by mitsubishi ufj asset management co ltd mitsubishi ufj asset management co ltd buy share of esco technology inc nyseese maryland state retirement pension system reduces position in toast inc nysetost maryland state retirement pension system reduces stock position in rollins inc nyserol mitsubishi ufj asset management co ltd raise stock holding in planet lab pbc nysepl mitsubishi ufj asset management co ltd grows stock holding in caretrust reit inc nasdaqctre mitsubishi ufj asset management co ltd reduces holding in xpeng inc sponsored adr nysexpev headtohead survey metallurgical corp of china mlluy the competition maryland state retirement pension system sell share of toast inc nysetost sun community inc nysesui stock position cut by maryland state retirement pension system maryland state retirement pension system cut stake in burlington store inc nyseburl maryland state retirement pension system sell share of nutanix nasdaqntnx reviewing metallurgical corp of china mlluy and it peer financial contrast nu skin enterprise nysenus and pitanium nasdaqptnm ehang holding limited unsponsored adr nasdaqeh receives average pt from analyst mitsubishi ufj asset management co ltd acquires share of asml holding nv nasdaqasml gold price hover near record high on fed rate cut expectation entain plc otcmktsgmvhf receives average rating of moderate buy from analyst h b fuller company nyseful vp sell in stock swiss re ltd otcmktsssrey receives consensus recommendation of hold from analyst erste group bank ag otcmktsebkdy receives consensus rating of moderate buy from brokerage svenska handelsbanken ab publ otcmktssvnly receives average recommendation of reduce from brokerage april underwood sell share of zillow group inc nasdaqz stock diana sue ferguson acquires share of sally beauty holding inc nysesbh stock jayshree ullal sell share of arista network inc nyseanet stock unity software inc nyseu coo sell in stock joao magalhaes sell share of h b fuller company nyseful stock anglo american otcmktsngloy receives average recommendation of hold from analyst svenska handelsbanken ab publ otcmktssvnly receives consensus rating of reduce from analyst diana sue ferguson acquires share of sally beauty holding inc nysesbh stock centerspace nysecsr ceo purchase in stock nordea bank ab otcmktsnrdby given consensus recommendation of moderate buy by analyst april underwood sell share of zillow group inc nasdaqz stock john f sheridan acquires share of tandem diabetes care inc nasdaqtndm stock idaho strategic resource inc nyseamericanidr director kevin g shiell sell share entain plc otcmktsgmvhf given average recommendation of moderate buy by brokerage alexander blum sell share of unity software inc nyseu stock wgn in chicago hire akemi harrison a news director lumen technology nyselumn grupo televisa nysetv head to head survey a hidden fuel source beneath the midwest scientist are investigating prediction this artificial intelligence ai stock could hit a trillion valuation by top wall street forecaster revamp brinker international expectation ahead of q earnings oi glass inc nyseoi receives average rating of moderate buy from analyst oi glass inc nyseoi given consensus rating of moderate buy by analyst zion bancorporation national association ut take position in charter communication inc nasdaqchtr zion bancorporation national association ut make new investment in mettlertoledo international inc nysemtd share in midamerica apartment community inc nysemaa purchased by zion bancorporation national association ut zion bancorporation national association ut acquires share of halliburton company nysehal zion bancorporation national association ut purchase share of everest group ltd nyseeg perplexity billion dream of buying chrome is never going to happen it also kind of brilliant oscar health inc nyseoscr share bought by connor clark lunn investment management ltd zion bancorporation national association ut take position in evergy inc nasdaqevrg zion bancorporation national association ut take position in garmin ltd nysegrmn connor clark lunn investment management ltd purchase share of commercecom inc nasdaqcmrc zion bancorporation national association ut acquires share of old dominion freight line inc nasdaqodfl zion bancorporation national association ut take position in halliburton company nysehal connor clark lunn investment management ltd increase position in tuya inc sponsored adr nysetuya zion bancorporation national association ut invests in eqt corporation nyseeqt cop cartel and the new drone war on the southern border chris murphy go all in on funding bill boycott a dems seek bipartisanship letter support is plummeting for war on palestine hamas terrorist posing a aid worker taken out by israel japan pokemon happy meal promo end poorly american eagle sydney sweeney ad divide young people poll federal policing take effect in washingtonwhat to know trump cant accept bad news here how that hurt the first amendment opinion nyt chief white house reporter challenged to walk dc at night after mocking trump today in history august east germany close berlin border country once ruled by cleopatra crossword clue partner on a political ticket crossword clue congress must hear from jeffrey epstein victim about ghislaine maxwell role opinion a america navy shrink time is running out to prevent communist china from seizing control of the pacific a fix is needed so candidate wont shun public financing for governor race steve collins presidential order on homelessness push u back in time opinion today in history august east germany close berlin border saving life no more rfk risk u all in targeting mrna vaccine research we have a responsibility to house mainer who live outside opinion letter posting commandment in school is part of learning column coal need to be part of america energy future column america need affordable energy coal doesnt fit the bill more than gop attorney general call on rfk jr fda to reinstate safeguard for abortion drug red state can now protect kid from trans procedure after court victory ap news summary at am edt rural homeowner share how he avoided a hour blackout without a gas generator no one else can do this pax trumpana trump nobel peace prize resume tennessee town approves deal to turn closed prison into immigration detention facility tennessee town approves deal to turn closed prison into immigration detention facility editorial money role in american politics could get even worse white house reviewing all smithsonian museum exhibit content before america celebration in anthony scaramucci share warning about maga going marxist and maoist yale professor alleges assault on freemarket capitalism morning brief low inflation tariff truce wchina dc takeover fallout plan for postwar gaza elizabeth warren say trump promised to slash power bill in half but instead cost rose aid program vanished and ceo cashed in zelenskyy to visit berlin for european virtual meeting ahead of the trumpputin summit german government say reader sound off on investigating tish james troop in dc and hamas culpability donald trump a strongman riling up his base and investigating his enemy out of our past august senegal prime minister move to increase tobacco tax to save life stimulate economic growth paul geisler woman might see themselves in basf planetarium show another viewpoint vaccine decision show rfk jr unfit jungbunzlauer awarded platinum medal from ecovadis for sustainability performance duffy toni marie ppr team preview helix highlander letter for wednesday august dea head say phone call are coming from all over the country from fed agent wanting to help in dc letter green energy is still needed just say no to more marijuana appeal court allows arkansas firstinthenation ban on gender transition care for minor to be enforced chisholm chisholm kilpatrick ltd and national veteran legal service program challenge va denial of benefit to caregiver after veteran death havard getting by with a little help from your friend trump at the kennedy center on the same day recipient of the honor are announced trump unveils nomination for judicial post california and trump spar in court over military deployment in la obituary lois eastman death notice shirley f bennett death notice clifford g foster lewis welcome to american politics without norm this small maine town had school budget trouble then the cost got higher trump official tell census worker congress ha final say over the count not trump san diego resident rally for killed al jazeera journalist dc crime forget issue scott jennings show dems are on the losing side of a one here why trump bls nominee ha ignited a cascade of criticism from economist trump administration threatens fine jail for washington homeless bernie sander say do this instead california governor demand trump abandon texas redistricting push threatens ballot measure response wife of south korea jailed expresident yoon arrested over corruption allegation leaving a top trump administration post the president may have an ambassadorship for you leaving a top trump administration post the president may have an ambassadorship for you judge to hear argument on halting alligator alcatraz construction over environmental concern judge to hear argument on halting alligator alcatraz construction over environmental concern judge to hear argument on halting alligator alcatraz construction over environmental concern judge to hear argument on halting alligator alcatraz construction over environmental concern alligator alcatraz migrant jail could halt construction over environmental concern construction at alligator alcatraz migrant jail could be halted over environmental concern construction at alligator alcatraz migrant jail could be halted over environmental concern judge to hear argument on halting alligator alcatraz construction over environmental concern judge to hear argument on halting alligator alcatraz construction over environmental concern charles gleason lawson ap news in brief at am edt clarence page time for democrat to get serious about their political future seth meyers give trump a failing score after leaving this name out of cognitive brag cartoon democrat are desperate for this letter nevada a sanctuary state trump is making major concession to union boss is it worth it letter gov newsom and the gender confusion club letter natural gas pay off for pa letter right age for sex letter where is our humanity cartoon dont trust him the word of the day letter union leader ignores psycho with hbombs tyler t ray caswell ha nh economy headed in the right direction congressman interacts with local official business leader letter social security now but can it long persist garrison keillor one favor lord if you have a moment editorial reviewjournal socialism lead to worse thing than expensive food can we somehow make backtoschool safe again israel emerging occupation consensus letter voting restriction will do more harm than good in nh len cannon more conservative point of view required fresh grocery arrive in san diego the same day with amazon cisco earnings are imminent these most accurate analyst revise forecast ahead of earnings call stock of the day new alltime high for amd advanced ceramic and nanoceramics lead the way in a market positioned for longterm growth amazon stock gain on big expansion of perishable delivery instacart grocery stock fall investor business daily saga communication inc declares a quarterly cash dividend of per share amazon expands it perishable delivery service putting pressure on traditional grocer ap news alibaba stock soar can ai innovation sustain the growth cold wallet m stage presale offer roi before price move higher in july saga communication inc declares a quarterly cash dividend of per share amazon expands it perishable delivery service putting pressure on traditional grocer tesla diner cut menu and hour amid surging demand in hollywood the nuke race is on energy department tap startup to build a reactor by mid barrons wait the u wind industry is coming back krg baghdad strike another oil deal podcastone post percent gain in q mario paglino and gianni grossi die designer made barbies into art worksport wksp q earnings call transcript apple plot expansion into ai robot home security and smart display bloombergcom south maintains lead in construction backlog driven by population growth new york sue zelle say security lapse led to billion consumer fraud loss reuters wescap management group inc ha position in johnson johnson nysejnj msh capital advisor llc acquires share of johnson johnson nysejnj hurricane owner tom dundon strike tentative deal to buy the portland trail blazer best thriller anime show on crunchyroll amazon roll out sameday perishable grocery delivery to more than u city nvidia face b sale threat a china target h chip over security inventhelp inventor develops new refueling accessory mho pure insurance flip the script on teen driving safetyfocused on education and safety not tracking amazon expands sameday grocery delivery rival stock plunge toll brother announces final opportunity to own a new luxury home at parklynn hill community in fountain inn south carolina connecticut green bank and goodleap partner to create an aipowered virtual power plant initiative to unlock saving for homeowner and enhance grid reliability toll brother announces final opportunity to own a new luxury home at parklynn hill community in fountain inn south carolina nissan cut ariya ev price by over with a new entrylevel model crypto exchange bullish valued at nearly billion in blowout nyse debut reuters deleteme v optery top tool to erase personal data from broker vietnam hightech push to become asia next tiger by leadingre add accounttech to solution group program princeton top linkedins first u college ranking for career rigetti computing stay the course in push for quantum advantage mobile momentum how u consumer and merchant are shaping unified commerce why share of globale online are sinking today well fargo company nysewfc share purchased by financial network wealth advisor llc debt surge past trillion a interest cost eclipse defense spending intellicheck report result tampa ma firm buy naming right to nhl arena the secure mat keep your bathroom clean and safe walmart offer permanent grocery discount to m u worker why cfo must stop treating compensation a a cost twin owner opt to halt the sale and keep the club in the family adding new investor instead ai drive layoff at microsoft and intel for efficiency twin owner opt to halt sale keep club in family public notice greenlight electricity centre project participant funding available costco stock is trading lower wednesday whats going on u housing market frozen by high rate awaits fed cut for thaw what the option market tell u about galaxy digital lowes companiess option a look at what the big money is thinking united airline holdingss option frenzy what you need to know looking into midamerica apartment recent short interest is the market bullish or bearish on insulet a look into borgwarner inc price over earnings what doe the market think about united parcel service best job and career in data science kazakhstan top central asia in digital advancement ecovyst report solid q tightens fullyear guidance new opportunity for business and cooperation are created with the visit to argentina of three central american and caribbean minister and the participation of international expert in a farmer conference way a chatbot will help your customer service team outdoor holding company announces preferred stock dividend top crypto coin to buy blockdag aave link and bnb setting the standard for community strength if you bought worth of shiba inu exactly year ago here what youd have one new coin promise far bigger return sp company see no recession a economy and earnings improve youtube add custom ctas to promote tool for higher ad conversion ap business summarybrief at pm edt this data center play offer fresh buy point but test a moving average boat trader unveils innovation to transform the boatbuying and selling experience u firm stockpile import at la port amid tariff fear att million settlement you could receive up to urban one report revenue drop in q playboy to relocate global headquarters to miami beach plan to develop new playboy club tampa businessman sentenced to federal prison for tax fraud casual dining take a hit stock to watch barrons overlooked stock that are soaring even more than ai save up to on the mindblowing galaxy book pro ai is joining cybersecurity team it could be a bad idea citrin cooperman expands west coast presence with acquisition of bpsd amazon expands sameday grocery delivery to u city secure investment management rank on the inc list of america fastestgrowing private company a leap in national ranking new york ag sue zelle parent over b fraud loss seek reform bolt stock is trading higher wednesday whats going on rigetti revenue drop percent in q little pepe lilpepe presale hit a investor snap up billion token to end stage early no shoe at the startup backtoschool deal alert save when you buy of school supply at amazon tom dundon to buy the portland trail blazer from paul allen estate report blazer find buyer staying in portland coffee sold in state recalled for potential glass fragment aolcom need a boost ankers back to school sale is slashing off charger power bank and more cima selected to lead a unique and major water research project bank of america corporation nysebac share purchased by mondrian investment partner ltd union pacific train derails in texas car off track no injury washington dc police union accuses city of cooking the book with crime stats kiss star say theyre grateful for kennedy center honor stevens democrat need a creative approach to the midterm election ltte frustrated by year of stalled street rotary and tree work fagan the heritage foundation founder legacy is complicated federal court allows doge access to sensitive data on million of american rfk jr is trying to make covid great again why did we build community that discourage walking the ethical life podcast the latest trump warns of very severe consequence if putin continues ukraine war australia punter playing at nebraska break down while discussing homesickness biden admin scrapped bestqualified standard for air traffic controller academy doc show flesheating bacteria found in swimmer who wa at cape cod beach trump shock troop coming to a city near you federal agent will be out on patrol in washington the white house say federal agent will be out on patrol in washington the white house say federal agent will be out on patrol in washington the white house say federal agent will be out on patrol in washington the white house say federal agent will be out on patrol in washington the white house say federal agent will be out on patrol in washington the white house say why democrat have limited power to fight republican redistricting trump real endgame with crime in dc is coming into view national guard presence will grow in dc wednesday night white house say child pregnant woman mistreated in immigration detention sen ossoff say what bangladesh ha achieved in the year since it revolution chesapeake council to allow member virtual participation in some meeting environmental concern could halt construction at florida alligator alcatraz immigration jail environmental concern could halt construction at florida alligator alcatraz immigration jail reporter share eyewitness account from front line of the russiaukraine war federal agent will be out on patrol in washington the white house say trump is engaged in stalinist attempt to rewrite u history axelrod warns live now texas house democrat hold news conference with indiana lawmaker trump visit kennedy center and announces honoree see the full list gaslighter king gavin newsom get ratioed into orbit for attacking stephen miller crime stat truth can public medium survive did india really shoot down an f fighting falcon fighter a freedom caucus target jackson hole housing policy gordon back community solution state leader take on record fee california wa a model for transparency now the capitol operates in the dark california wa a model for transparency now the capitol operates in the dark port a deep dive into how the sausage get made in north dakota trump say he doesnt believe he can convince putin to stop bombing ukrainian civilian att million settlement you could receive up to tulsi gabbard release damning email that put clapper in the hot seat judge is skeptical of justice department lawsuit against federal judge a trump try to limit power of judiciary trump threatens very severe consequence if russia doesnt agree to end ukraine war trump predicts little progress in potential shutdown talk with crazy schumer jeffries zohran mamdani take aim at trump immigration policy in fiery staten island speech we are fighting to keep nyc a sanctuary city auditor california could save million a year by letting state employee work remote trump warns putin of very severe consequence if ukraine war continues trump pledged to move homeless people from washington what we know and dont know about his plan trump pledged to move homeless people from washington what we know and dont know about his plan trump pledged to move homeless people from washington what we know and dont know about his plan trump pledged to move homeless people from washington what we know and dont know about his plan trump pledged to move homeless people from washington what we know and dont know about his plan trump pledged to move homeless people from washington what we know and dont know about his plan trump pledged to move homeless people from washington what we know and dont know about his plan trump pledged to move homeless people from washington what we know and dont know about his plan trump pledged to move homeless people from washington what we know and dont know about his plan new york claim zelles shoddy security enabled a billion dollar in scam wall street journal give trump a stern warning new connection credit union clubhouse coming to riverside golf course judge skeptical of trump administration suit against maryland judge will rule by labor day judge skeptical of trump administration suit against maryland judge will rule by labor day judge skeptical of trump administration suit against maryland judge will rule by labor day judge skeptical of trump administration suit against maryland judge will rule by labor day judge skeptical of trump administration suit against maryland judge will rule by labor day judge skeptical of trump administration suit against maryland judge will rule by labor day he a bit of a mouth local react to vance visit president trump announces sylvester stallone kiss a kennedy center honoree trump name kennedy center honoree say hell host award show brit throw vance a not welcome party on his latest vacation trump tariff are forcing canada to address it money laundering problem trump administration lawsuit against all of maryland federal judge meet skepticism in court u senate approves m to expand willmar veterinary testing lab trump pledged to move homeless people from washington what we know and dont know about his plan allegedly nonpartisan common cause unmasked cheering leftist gerrymandering while condemning the gop pocatello model railroad and historical society open house set for saturday ap news summary at pm edt the latest trump warns of very severe consequence if putin continues ukraine war nyc ramp up effort to remove mentally ill from street subway trumpputin summit on ukraine is latest chapter in alaska long history and tension with russia trumpputin summit on ukraine is latest chapter in alaska long history and tension with russia trumpputin summit on ukraine is latest chapter in alaska long history and tension with russia trumpputin summit on ukraine is latest chapter in alaska long history and tension with russia trumpputin summit on ukraine is latest chapter in alaska long history and tension with russia trumpputin summit on ukraine is latest chapter in alaska long history and tension with russia trumpputin summit on ukraine is latest chapter in alaska long history and tension with russia trumpputin summit on ukraine is latest chapter in alaska long history and tension with russia trumpputin summit on ukraine is latest chapter in alaska long history and tension with russia epa terminates center for rural affair million solar grant trump look to extend dc police takeover beyond day trump look to extend dc police takeover beyond day harvard and the trump administration are nearing a settlement including a million payment harvard and the trump administration are nearing a settlement including a million payment harvard and the trump administration are nearing a settlement including a million payment trump detail expectation for meeting with putin chicago rapper king yella hold chest fall over after drinking crown royal throughout interview a trump turn dc into a fascist laboratory we must mount our own defense what to know about the putintrump summit in alaska maga hang mission accomplished hype on trump tariff jonah goldberg
'@

$ws.Cells.Item($row, 1).Value = 45883
$ws.Cells.Item($row, 2).Value = 6468.5400390625
$ws.Cells.Item($row, 3).Value = 6473.919921875
$ws.Cells.Item($row, 4).Value = 6441.06982421875
$ws.Cells.Item($row, 5).Value = 6453.4599609375
$ws.Cells.Item($row, 6).Value = 2766449000
$ws.Cells.Item($row, 7).Value = 0.0003030908012924
$ws.Cells.Item($row, 8).Value = 1
$ws.Cells.Item($row, 11).Value = 0.123
$ws.Cells.Item($row, 12).Value = 0.095
$ws.Cells.Item($row, 13).Value = 0.782
$ws.Cells.Item($row, 14).Value = 0.998
$ws.Cells.Item($row, 15).Value = 129
$ws.Cells.Item($row, 16).Value = 123
$ws.Cells.Item($row, 17).Value = 0
$ws.Cells.Item($row, 18).Value = 315
$ws.Cells.Item($row, 19).Value = 99
$ws.Cells.Item($row, 20).Value = 424
$ws.Cells.Item($row, 21).Value = 193
$ws.Cells.Item($row, 22).Value = 98
$ws.Cells.Item($row, 23).Value = 40
$ws.Cells.Item($row, 24).Value = 145

# Preserve the date number format on column A like the preceding rows
$ws.Range("A37").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Headline columns (I and J) share the same long text blob
$ws.Cells.Item($row, 9).Value = $headline
$ws.Cells.Item($row, 10).Value = $headline

Write-Host "Row 37 written"